$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-24 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-25 Saturday", 2) | Out-Null
$d.Content.Find.Execute("970÷7=138, 4", $true, $false, $false, $false, $false, $true, 1, $false, "208÷8=26, 0", 2) | Out-Null
$d.Content.Find.Execute("535÷6=89, 1", $true, $false, $false, $false, $false, $true, 1, $false, "244÷2=122, 0", 2) | Out-Null
$d.Content.Find.Execute("432÷8=54, 0", $true, $false, $false, $false, $false, $true, 1, $false, "897÷3=299, 0", 2) | Out-Null
$d.Content.Find.Execute("256÷8=32, 0", $true, $false, $false, $false, $false, $true, 1, $false, "420÷7=60, 0", 2) | Out-Null
$d.Content.Find.Execute("984÷5=196, 4", $true, $false, $false, $false, $false, $true, 1, $false, "502÷4=125, 2", 2) | Out-Null
$d.Content.Find.Execute("675÷8=84, 3", $true, $false, $false, $false, $false, $true, 1, $false, "530÷2=265, 0", 2) | Out-Null
$d.Content.Find.Execute("328÷2=164, 0", $true, $false, $false, $false, $false, $true, 1, $false, "805÷9=89, 4", 2) | Out-Null
$d.Content.Find.Execute("346÷6=57, 4", $true, $false, $false, $false, $false, $true, 1, $false, "830÷2=415, 0", 2) | Out-Null
$d.Content.Find.Execute("196÷8=24, 4", $true, $false, $false, $false, $false, $true, 1, $false, "341÷4=85, 1", 2) | Out-Null
$d.Content.Find.Execute("917÷9=101, 8", $true, $false, $false, $false, $false, $true, 1, $false, "769÷3=256, 1", 2) | Out-Null
$d.Content.Find.Execute("394÷9=43, 7", $true, $false, $false, $false, $false, $true, 1, $false, "270÷6=45, 0", 2) | Out-Null
$d.Content.Find.Execute("524÷8=65, 4", $true, $false, $false, $false, $false, $true, 1, $false, "205÷4=51, 1", 2) | Out-Null
$d.Content.Find.Execute("266÷3=88, 2", $true, $false, $false, $false, $false, $true, 1, $false, "458÷3=152, 2", 2) | Out-Null
$d.Content.Find.Execute("628÷4=157, 0", $true, $false, $false, $false, $false, $true, 1, $false, "293÷3=97, 2", 2) | Out-Null
$d.Content.Find.Execute("783÷4=195, 3", $true, $false, $false, $false, $false, $true, 1, $false, "525÷8=65, 5", 2) | Out-Null
$d.Content.Find.Execute("125÷9=13, 8", $true, $false, $false, $false, $false, $true, 1, $false, "628÷5=125, 3", 2) | Out-Null
$d.Content.Find.Execute("739÷2=369, 1", $true, $false, $false, $false, $false, $true, 1, $false, "679÷8=84, 7", 2) | Out-Null
$d.Content.Find.Execute("670÷6=111, 4", $true, $false, $false, $false, $false, $true, 1, $false, "187÷7=26, 5", 2) | Out-Null
$d.Content.Find.Execute("917÷4=229, 1", $true, $false, $false, $false, $false, $true, 1, $false, "502÷7=71, 5", 2) | Out-Null
$d.Content.Find.Execute("708÷8=88, 4", $true, $false, $false, $false, $false, $true, 1, $false, "712÷6=118, 4", 2) | Out-Null
$d.Content.Find.Execute("969÷3=323, 0", $true, $false, $false, $false, $false, $true, 1, $false, "105÷7=15, 0", 2) | Out-Null
$d.Content.Find.Execute("497÷6=82, 5", $true, $false, $false, $false, $false, $true, 1, $false, "173÷6=28, 5", 2) | Out-Null
$d.Content.Find.Execute("701÷2=350, 1", $true, $false, $false, $false, $false, $true, 1, $false, "665÷6=110, 5", 2) | Out-Null
$d.Content.Find.Execute("212÷2=106, 0", $true, $false, $false, $false, $false, $true, 1, $false, "665÷4=166, 1", 2) | Out-Null
$d.Content.Find.Execute("857÷6=142, 5", $true, $false, $false, $false, $false, $true, 1, $false, "154÷6=25, 4", 2) | Out-Null
